$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.585.34"
$ws.Range("E2").Value = "  -4.16%  "

$ws.Range("D3").Value = "3.315.85"
$ws.Range("E3").Value = "  -1.04%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "572.93"
$ws.Range("E5").Value = "  -3.74%  "

$ws.Range("D6").Value = "182.08"
$ws.Range("E6").Value = "  -5.61%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("D9").Value = "3.315.55"
$ws.Range("E9").Value = "  -1.01%  "

$ws.Range("E10").Value = "  -2.95%  "

$ws.Range("D11").Value = "6.64"
$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("E12").Value = "  -5.33%  "

$ws.Range("D13").Value = "3.895.63"
$ws.Range("E13").Value = "  -0.94%  "

$ws.Range("D14").Value = "0.138"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "27.08"
$ws.Range("E15").Value = "  -4.68%  "

$ws.Range("D16").Value = "66.662.46"
$ws.Range("E16").Value = "  -4.05%  "

$ws.Range("E17").Value = "  -2.66%  "

$ws.Range("D18").Value = "3.304.90"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").Value = "439.84"
$ws.Range("E19").Value = "  +2.68%  "

$ws.Range("D20").Value = "13.69"
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("E21").Value = "  -2.91%  "

$ws.Range("D22").Value = "7.64"
$ws.Range("E22").Value = "  -1.31%  "

$ws.Range("D23").Value = "73.82"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").Value = "3.458.33"
$ws.Range("E25").Value = "  -1.03%  "

$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("D28").Value = "0.192"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").Value = "9.01"
$ws.Range("E29").Value = "  -6.37%  "

$ws.Range("E30").Value = "  -1.02%  "

$ws.Range("E31").Value = "  -2.70%  "

$ws.Range("E32").Value = "  -0.98%  "

$ws.Range("D33").Value = "5.33"
$ws.Range("E33").Value = "  -5.88%  "

$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").Value = "6.77"
$ws.Range("E35").Value = "  -4.11%  "

$ws.Range("E36").Value = "  -5.03%  "

$ws.Range("E37").Value = "  -1.09%  "

$ws.Range("D38").Value = "160.41"
$ws.Range("E38").Value = "  -2.45%  "

$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").Value = "1.84"
$ws.Range("E40").Value = "  -4.43%  "

$ws.Range("D41").Value = "2.809.11"
$ws.Range("E41").Value = "  +2.24%  "

$ws.Range("D42").Value = "0.789"
$ws.Range("E42").Value = "  -2.76%  "

$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("D44").Value = "6.24"
$ws.Range("E44").Value = "  -3.73%  "

$ws.Range("E45").Value = "  -1.92%  "

$ws.Range("D46").Value = "40.14"
$ws.Range("E46").Value = "  -2.54%  "

$ws.Range("D47").Value = "24.21"
$ws.Range("E47").Value = "  -5.24%  "

$ws.Range("E48").Value = "  -7.07%  "

$ws.Range("D49").Value = "318.33"
$ws.Range("E49").Value = "  -7.38%  "

$ws.Range("D50").Value = "0.0272"
$ws.Range("E50").Value = "  -3.74%  "

$ws.Range("D51").Value = "0.982"
$ws.Range("E51").Value = "  -3.02%  "
